$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B. This shifts the existing
# "PercActivations" column (B) to C and "PercSegmentAreas" column (C) to D.
$ws.Columns("B").Insert()

# Column A currently still holds the segment-name text (with its header
# style carried over). Move those names into the freshly inserted column B,
# and replace column A's contents with the numeric segment index (0-based),
# keeping A's existing header/border style.
for ($r = 2; $r -le 20; $r++) {
  $name = $ws.Cells.Item($r, 1).Value()
  $ws.Cells.Item($r, 2).Value = $name
  # The Insert() operation copied column A's style onto the new column B;
  # the data values in B should be plain (unstyled), so clear that back out.
  $ws.Cells.Item($r, 2).ClearFormats()
  $ws.Cells.Item($r, 1).Value = $r - 2
}

# Give the new column its header label, matching the style already used by
# the other header cells in row 1.
$ws.Cells.Item(1, 2).Value = "segments"
$ws.Range("C1").Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
